$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    # Force the cell to be treated as literal text so numeric-looking
    # strings (e.g. "357.28") aren't silently coerced into floating point
    # numbers (which would round-trip as "357.27999999999997" etc.), and
    # restore the default ("Normal") style afterwards so we don't leave a
    # stray number-format style attached to the cell.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "52.286.92"
$ws.Range("E2").Value = "  +1.39%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "2.832.50"
$ws.Range("E3").Value = "  +3.28%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.01%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "357.28"
$ws.Range("E5").Value = "  +7.62%  "

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "114.57"
$ws.Range("E6").Value = "  -1.71%  "

# Row 7 - XRP
Set-TextValue $ws.Range("D7") "0.548"
$ws.Range("E7").Value = "  +2.89%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.06%  "

# Row 9 - Cardano
Set-TextValue $ws.Range("D9") "0.605"
$ws.Range("E9").Value = "  +5.72%  "

# Row 10 - Avalanche
Set-TextValue $ws.Range("D10") "42.04"
$ws.Range("E10").Value = "  +1.68%  "

# Row 11 - Dogecoin
Set-TextValue $ws.Range("D11") "0.0849"
$ws.Range("E11").Value = "  +2.09%  "

# Row 12 - was TRON, now Chainlink
$ws.Range("B12").Value = "Chainlink"
$ws.Range("C12").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws.Range("D12") "19.96"
$ws.Range("E12").Value = "  -1.10%  "

# Row 13 - was Chainlink, now TRON
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue $ws.Range("D13") "0.131"
$ws.Range("E13").Value = "  +1.56%  "

# Row 14 - Polkadot
Set-TextValue $ws.Range("D14") "7.82"
$ws.Range("E14").Value = "  +3.54%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D15") "3.269.68"
$ws.Range("E15").Value = "  +3.10%  "

# Row 16 - WrappedEther
Set-TextValue $ws.Range("D16") "2.832.82"
$ws.Range("E16").Value = "  +3.52%  "

# Row 17 - Polygon
Set-TextValue $ws.Range("D17") "0.899"
$ws.Range("E17").Value = "  +1.89%  "

# Row 18 - WrappedBTC
Set-TextValue $ws.Range("D18") "52.050.28"
$ws.Range("E18").Value = "  +1.11%  "

# Row 19 - ImmutableX
Set-TextValue $ws.Range("D19") "3.18"
$ws.Range("E19").Value = "  +2.46%  "

# Row 20 - Uniswap
Set-TextValue $ws.Range("D20") "7.31"
$ws.Range("E20").Value = "  +7.31%  "

# Row 21 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D21") "13.74"
$ws.Range("E21").Value = "  +2.42%  "

# Row 22 - ShibaInu
Set-TextValue $ws.Range("D22") "0.0₃0990"
$ws.Range("E22").Value = "  +3.02%  "

# Row 23 - BitcoinCash
Set-TextValue $ws.Range("D23") "270.30"
$ws.Range("E23").Value = "  -2.48%  "

# Row 24 - Litecoin
Set-TextValue $ws.Range("D24") "69.66"
$ws.Range("E24").Value = "  +0.36%  "

# Row 25 - PancakeSwap
Set-TextValue $ws.Range("D25") "2.81"
$ws.Range("E25").Value = "  +6.53%  "

# Row 26 - EthereumClassic
Set-TextValue $ws.Range("D26") "26.89"
$ws.Range("E26").Value = "  +0.77%  "

# Row 27 - Dai
Set-TextValue $ws.Range("D27") "1.00"
$ws.Range("E27").Value = "  +0.12%  "

# Row 28 - Cosmos
Set-TextValue $ws.Range("D28") "10.25"
$ws.Range("E28").Value = "  +0.87%  "

# Row 29 - Toncoin
Set-TextValue $ws.Range("D29") "2.25"
$ws.Range("E29").Value = "  +1.37%  "

# Row 30 - Kaspa
$ws.Range("E30").Value = "  +0.66%  "

# Row 31 - OKB
$ws.Range("E31").Value = "  +0.95%  "

# Row 32 - InjectiveProtocol
Set-TextValue $ws.Range("D32") "33.96"
$ws.Range("E32").Value = "  -2.79%  "

# Row 33 - Filecoin
Set-TextValue $ws.Range("D33") "5.88"
$ws.Range("E33").Value = "  +6.03%  "

# Row 34 - VeChain
Set-TextValue $ws.Range("D34") "0.0437"
$ws.Range("E34").Value = "  +26.23%  "

# Row 35 - Hedera
Set-TextValue $ws.Range("D35") "0.0830"
$ws.Range("E35").Value = "  +1.21%  "

# Row 36 - FirstDigitalUSD
Set-TextValue $ws.Range("D36") "0.999"
$ws.Range("E36").Value = "  -0.21%  "

# Row 37 - ARBITRUM
Set-TextValue $ws.Range("D37") "2.10"
$ws.Range("E37").Value = "  +1.27%  "

# Row 38 - RenderToken
Set-TextValue $ws.Range("D38") "4.90"
$ws.Range("E38").Value = "  +0.16%  "

# Row 39 - Celestia
Set-TextValue $ws.Range("D39") "18.52"
$ws.Range("E39").Value = "  -2.60%  "

# Row 40 - LidoDAOToken
$ws.Range("E40").Value = "  +1.50%  "

# Row 41 - Stacks
Set-TextValue $ws.Range("D41") "2.58"
$ws.Range("E41").Value = "  +7.42%  "

# Row 42 - EnergySwap
Set-TextValue $ws.Range("D42") "23.49"
$ws.Range("E42").Value = "  +1.03%  "

# Row 43 - Stellar
Set-TextValue $ws.Range("D43") "0.116"
$ws.Range("E43").Value = "  +2.16%  "

# Row 44 - Monero
Set-TextValue $ws.Range("D44") "126.86"
$ws.Range("E44").Value = "  -1.71%  "

# Row 45 - WEMIXToken
Set-TextValue $ws.Range("D45") "2.29"
$ws.Range("E45").Value = "  +2.12%  "

# Row 46 - NEARProtocol
Set-TextValue $ws.Range("D46") "3.37"
$ws.Range("E46").Value = "  +1.25%  "

# Row 47 - Maker
Set-TextValue $ws.Range("D47") "2.048.13"
$ws.Range("E47").Value = "  -2.48%  "

# Row 48 - ApeXProtocol
$ws.Range("E48").Value = "  +3.77%  "

# Row 49 - SEI
Set-TextValue $ws.Range("D49") "0.958"
$ws.Range("E49").Value = "  +10.87%  "

# Row 50 - THORChain
Set-TextValue $ws.Range("D50") "5.73"
$ws.Range("E50").Value = "  +3.99%  "

# Row 51 - FraxShare
$ws.Range("E51").Value = "  +0.02%  "
